# Apply entropy calculation results update ("all entropy calculations have been performed")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E13 tiny precision update
$ws.Range("E13").Value = [double]"3.586954457787532e-12"

# Row 16 - alpha (EEG1, session 2): results now computed
$ws.Range("E16").Value = [double]"0.7450723481033423"
$ws.Range("F16").Value = [double]"-0.4333408326701519"
$ws.Range("G16").Value = [double]"79.97604888891956"

# Row 17 - beta (EEG1, session 2): results now computed, now significant
$ws.Range("E17").Value = [double]"0.000643871801222988"
$ws.Range("F17").Value = [double]"-5.642684565960459"
$ws.Range("G17").Value = [double]"51.22433748035346"
$ws.Range("H17").Value = $true

# Row 18 - gamma (EEG1, session 2): results now computed, now significant
$ws.Range("E18").Value = [double]"0.02241830351434788"
$ws.Range("F18").Value = [double]"-3.942617668930942"
$ws.Range("G18").Value = [double]"39.29620486710468"
$ws.Range("H18").Value = $true

# Row 19 - all bands (EEG1, session 2): results now computed, now significant
$ws.Range("E19").Value = [double]"1.586969268842816e-15"
$ws.Range("F19").Value = [double]"-37.62292601181118"
$ws.Range("G19").Value = [double]"106.3364803441259"
$ws.Range("H19").Value = $true

# Row 22 - alpha (EEG2, session 2): results now computed, now significant
$ws.Range("E22").Value = [double]"6.461284622464442e-05"
$ws.Range("F22").Value = [double]"-4.531944186563699"
$ws.Range("G22").Value = [double]"57.32272466310992"
$ws.Range("H22").Value = $true

# E23 tiny precision update
$ws.Range("E23").Value = [double]"0.005181580316172342"

# E24 tiny precision update
$ws.Range("E24").Value = [double]"0.2446778414505822"

# E29 tiny precision update
$ws.Range("E29").Value = [double]"2.658171682420209e-13"

# F34 tiny precision update
$ws.Range("F34").Value = [double]"-6.029493151688476"

# E35 tiny precision update
$ws.Range("E35").Value = [double]"7.242078962700626e-49"
